$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 3.3
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 2.05
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 2.75
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62
$ws.Range("S2").Value = 4
$ws.Range("T2").Value = 1.22
$ws.Range("U2").Value = 1.5
$ws.Range("V2").Value = 2.5
$ws.Range("W2").Value = 1.95
$ws.Range("X2").Value = 1.8
$ws.Range("Y2").Value = 8.5
$ws.Range("AA2").Value = 12
$ws.Range("AC2").Value = 29
$ws.Range("AD2").Value = 41
$ws.Range("AE2").Value = 8.5
$ws.Range("AG2").Value = 17
$ws.Range("AI2").Value = 351
$ws.Range("AJ2").Value = 7
$ws.Range("AK2").Value = 10
$ws.Range("AN2").Value = 21
$ws.Range("AO2").Value = 34
